# Applies the PlatformDatabase2 update described in the commit
# "Addressing issue #21. Hope to have fixed it"
#
# Summary of content changes:
#  - Platform sheet: "Disaster Effect" (col J) values updated for several
#    platform rows (4,9 -> 9,13 ; 4 -> 13 ; Motorcycle gains a leading "1,").
#  - Platform sheet: a brand-new platform row ("Human", id 15) is appended
#    below Motorcycle.
#  - View/selection bookkeeping (zoom + active cell) on all three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Platform"
# ---------------------------------------------------------------------
$platform = $wb.Worksheets.Item("Platform")
$platform.Activate()

# -- Update "Disaster Effect" (column J) values for existing rows -------
$platform.Range("J2").Value = "9,13"
$platform.Range("J3").Value = "9,13"
$platform.Range("J4").Value = 13
$platform.Range("J5").Value = "9,13"
$platform.Range("J6").Value = 13
$platform.Range("J7").Value = "9,13"
$platform.Range("J8").Value = 13
$platform.Range("J9").Value = 13

# -- Append the new "Human" platform row (row 16) ------------------------
# Pull the formatting for the new row from the row above (Motorcycle, row
# 15) so fills/alignment match the rest of the table.
$platform.Range("A15:D15").Copy()
$platform.Range("A16:E16").PasteSpecial(-4122)
$platform.Range("G15:H15").Copy()
$platform.Range("G16:H16").PasteSpecial(-4122)
$platform.Range("J15").Copy()
$platform.Range("J16").PasteSpecial(-4122)
$platform.Range("L15").Copy()
$platform.Range("L16").PasteSpecial(-4122)
$platform.Range("M15").Copy()
$platform.Range("M16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$platform.Range("A16").Value = 15
$platform.Range("B16").Value = "Human"
$platform.Range("C16").Value = 1
$platform.Range("D16").Value = 20
$platform.Range("E16").Value = "G"
$platform.Range("G16").Value = "9,13"
$platform.Range("J16").Value = "1,2,10,11,12"
$platform.Range("H16").Value = "7,13"
$platform.Range("L16").Value = 10
$platform.Range("M16").Value = 8

# Update the "Motorcycle" row's Disaster Effect last, since this string is
# introduced into the shared-string table after the new row's strings.
$platform.Range("J15").Value = "1,2,5,10,11"

# -- View/selection bookkeeping ------------------------------------------
$platform.Range("H18").Select()
$excel.ActiveWindow.Zoom = 160

# ---------------------------------------------------------------------
# Sheet 2: "Custom Column Description"
# ---------------------------------------------------------------------
$customCol = $wb.Worksheets.Item("Custom Column Description")
$customCol.Activate()
$customCol.Range("B16").Select()

# ---------------------------------------------------------------------
# Sheet 3: "Supported Data Type"
# ---------------------------------------------------------------------
$dataType = $wb.Worksheets.Item("Supported Data Type")
$dataType.Activate()
$dataType.Range("C4").Select()

# ---------------------------------------------------------------------
# Leave the "Platform" sheet as the active / selected sheet, matching the
# tabSelected flag in the target workbook.
# ---------------------------------------------------------------------
$platform.Activate()
